# Fix the algorithm/conditions on filtering the status of candidates (1/12/2025)
# Rewrites the candidate status table (rows 3-21) to the corrected data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: status updated in place
$ws.Range("E3").Value = "3rd Interview"

# Rows 4-21: full refresh of Job ID / Company / Job Title / Candidate / Status / Action Date
$data = New-Object 'object[,]' 18,6

# Row 4
$data[0,0] = 681
$data[0,1] = "Metaview"
$data[0,2] = "Metaview - UK MM / Enterprise AE"
$data[0,3] = "Federico  Belluci"
$data[0,4] = "4th Interview"
$data[0,5] = 45989

# Row 5
$data[1,0] = 687
$data[1,1] = "CodeRabbit"
$data[1,2] = "Commercial AE x5 Bay Area"
$data[1,3] = "Ben Sawyer"
$data[1,4] = "3rd Interview"
$data[1,5] = 45978

# Row 6
$data[2,0] = 687
$data[2,1] = "CodeRabbit"
$data[2,2] = "Commercial AE x5 Bay Area"
$data[2,3] = "Jackson  Smith"
$data[2,4] = "CV Sent"
$data[2,5] = 45987

# Row 7
$data[3,0] = 687
$data[3,1] = "CodeRabbit"
$data[3,2] = "Commercial AE x5 Bay Area"
$data[3,3] = "Kyle Brown"
$data[3,4] = "1st Interview"
$data[3,5] = 45995

# Row 8
$data[4,0] = 721
$data[4,1] = "Rox"
$data[4,2] = "ENT AE East"
$data[4,3] = "Graham Hill"
$data[4,4] = "CV Sent"
$data[4,5] = 45973

# Row 9
$data[5,0] = 731
$data[5,1] = "Oscilar"
$data[5,2] = "Enterprise AE x5"
$data[5,3] = "Ryan Finkelstein"
$data[5,4] = "1st Interview"
$data[5,5] = 45992

# Row 10
$data[6,0] = 740
$data[6,1] = "Axion Ray"
$data[6,2] = "Enterprise Account Executive (East)"
$data[6,3] = "Mustafa Hubaishi"
$data[6,4] = "CV Sent"
$data[6,5] = 45985

# Row 11
$data[7,0] = 740
$data[7,1] = "Axion Ray"
$data[7,2] = "Enterprise Account Executive (East)"
$data[7,3] = "Graham Hill"
$data[7,4] = "CV Sent"
$data[7,5] = 45987

# Row 12
$data[8,0] = 751
$data[8,1] = "Materialize"
$data[8,2] = "CS3 Materialize - Enterprise Account Executive x3"
$data[8,3] = "Ryan Finkelstein"
$data[8,4] = "1st Interview"
$data[8,5] = 45993

# Row 13
$data[9,0] = 766
$data[9,1] = "Cogent Security"
$data[9,2] = "Enterprise Account Executive (US)"
$data[9,3] = "Daniel Smith"
$data[9,4] = "2nd Interview"
$data[9,5] = 45992

# Row 14
$data[10,0] = 776
$data[10,1] = "PointFive"
$data[10,2] = "Enterprise AE (EST)"
$data[10,3] = "Daniel Smith"
$data[10,4] = "3rd Interview"
$data[10,5] = 45992

# Row 15
$data[11,0] = 783
$data[11,1] = "Port"
$data[11,2] = "Mid-Market AE"
$data[11,3] = "Christopher Blair"
$data[11,4] = "1st Interview"
$data[11,5] = 45980

# Row 16
$data[12,0] = 785
$data[12,1] = "Dash0"
$data[12,2] = "Commercial AE Amsterdam"
$data[12,3] = "Oskar Stepien"
$data[12,4] = "4th Interview"
$data[12,5] = 45982

# Row 17
$data[13,0] = 792
$data[13,1] = "Energy Robotics"
$data[13,2] = "Partner Manager (UAE)"
$data[13,3] = "Hashem Abughazaleh"
$data[13,4] = "3rd Interview"
$data[13,5] = 45989

# Row 18
$data[14,0] = 796
$data[14,1] = "Redwood Software"
$data[14,2] = "Enterprise AE UK (Finance Automation)"
$data[14,3] = "Greg Keene"
$data[14,4] = "CV Sent"
$data[14,5] = 45987

# Row 19
$data[15,0] = 808
$data[15,1] = "Energy Robotics"
$data[15,2] = "Partner Manager (Houston)"
$data[15,3] = "Ahmed Mirza"
$data[15,4] = "3rd Interview"
$data[15,5] = 45993

# Row 20
$data[16,0] = 808
$data[16,1] = "Energy Robotics"
$data[16,2] = "Partner Manager (Houston)"
$data[16,3] = "Shashwat Anand"
$data[16,4] = "1st Interview"
$data[16,5] = 45993

# Row 21
$data[17,0] = 808
$data[17,1] = "Energy Robotics"
$data[17,2] = "Partner Manager (Houston)"
$data[17,3] = "[deleted]"
$data[17,4] = "1st Interview"
$data[17,5] = 45986

$ws.Range("A4:F21").Value = $data

# New rows 19-21 need the same Action Date number format as the rest of column F
$ws.Range("F18").Copy()
$ws.Range("F19:F21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("F19").Value = 45993
$ws.Range("F20").Value = 45993
$ws.Range("F21").Value = 45986

$wb.Save()
